# Update crypto price/volume snapshot values (GitHub Actions refresh).
# D-column values are stored as leading-apostrophe text (Excel "quote prefix")
# so numeric-looking strings like "21.83" or "1.001" stay text instead of
# being auto-converted to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.001.93'
$ws.Range("E2").Value = '  -1.88%  '
$ws.Range("D3").Value = '''1.830.23'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '''324.12'
$ws.Range("E5").Value = '  -2.84%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.3864'
$ws.Range("E8").Value = '  -1.43%  '
$ws.Range("D9").Value = '''0.07871'
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").Value = '''0.9592'
$ws.Range("D11").Value = '''21.83'
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").Value = '''1.793.48'
$ws.Range("E12").Value = '  -6.79%  '
$ws.Range("D13").Value = '''5.668'
$ws.Range("E13").Value = '  -3.09%  '
$ws.Range("D14").Value = '''6.897'
$ws.Range("D15").Value = '''0.06837'
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").Value = '''87.29'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D18").Value = '''0.000009918'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("E19").Value = '  -2.64%  '
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '''28.002.81'
$ws.Range("E21").Value = '  -2.02%  '
$ws.Range("D22").Value = '''5.313'
$ws.Range("E22").Value = '  -1.38%  '
$ws.Range("D23").Value = '''10.97'
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("D25").Value = '''2.038.16'
$ws.Range("E25").Value = '  -7.37%  '
$ws.Range("D26").Value = '''153.91'
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = '''19.11'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").Value = '''5.711'
$ws.Range("E28").Value = '  -6.64%  '
$ws.Range("D29").Value = '''1.957'
$ws.Range("E29").Value = '  -2.81%  '
$ws.Range("D30").Value = '''117.69'
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("D31").Value = '''0.09261'
$ws.Range("E31").Value = '  -1.72%  '
$ws.Range("D32").Value = '''0.9341'
$ws.Range("E32").Value = '  -4.14%  '
$ws.Range("D33").Value = '''5.274'
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("D34").Value = '''1.314'
$ws.Range("E34").Value = '  -2.22%  '
$ws.Range("D35").Value = '''3.294'
$ws.Range("E35").Value = '  -5.97%  '
$ws.Range("D36").Value = '''0.05861'
$ws.Range("E36").Value = '  -4.16%  '
$ws.Range("D37").Value = '''0.02140'
$ws.Range("E37").Value = '  -2.46%  '
$ws.Range("D38").Value = '''1.141'
$ws.Range("E38").Value = '  -1.90%  '
$ws.Range("D39").Value = '''7.776'
$ws.Range("E39").Value = '  +2.39%  '
$ws.Range("D40").Value = '''0.5577'
$ws.Range("D41").Value = '''9.846'
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("E42").Value = '  -1.79%  '
$ws.Range("D43").Value = '''11.62'
$ws.Range("E43").Value = '  -1.43%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = '''0.07022'
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.5252'
$ws.Range("E45").Value = '  -2.35%  '
$ws.Range("D46").Value = '''2.116'
$ws.Range("E46").Value = '  -11.22%  '
$ws.Range("D47").Value = '''1.822'
$ws.Range("E47").Value = '  -4.27%  '
$ws.Range("D48").Value = '''112.93'
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("D49").Value = '''1.103'
$ws.Range("E49").Value = '  -11.95%  '
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").Value = '''2.322'
$ws.Range("E51").Value = '  +0.38%  '
